$d = $word.ActiveDocument

# Locate the paragraph whose text is exactly "2." (the ticket-number heading)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "2.`r") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs.First
}

$r = $target.Range
# Exclude the trailing paragraph mark from the range so InsertXML only
# replaces the run content, leaving the paragraph (and its pPr) intact.
$r.MoveEnd(1, -1)

$newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>2</w:t></w:r>
            <w:r><w:t>2</w:t></w:r>
            <w:r><w:t>.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> Философский иррационализм А. Шопенгауэра</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r.InsertXML($newXml)
